# Automatic map update (2025-08-21 08:06:57)
# - Rows 3-8 get their data re-synced (re-ordered) from the upstream source.
# - Two brand-new rows (69, 70) are appended.
# - Sheet dimension grows from A1:P68 to A1:P70 automatically as the new
#   rows are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for every touched row (1-based sheet rows).
$rowsData = @(
    @{ Row=3;  A="5589";      B="12/31/2023"; C="ARCOS 1520";                     D="13"; E="799540526"; F="NEW"; G="Pendiente de Traspaso PROPIO"; H="Picada";          I=0; J="Cambio"; K="Nodo Teco";    L="Pasante"; M=-58.449125; N=-34.565958; O="Colegiales"; P="Capital Norte" }
    @{ Row=4;  A="4862";      B="1/23/2025";  C="ARCOS 2263";                     D="13"; E="802857379"; F="NEW"; G="Pendiente de Traspaso PROPIO"; H="picada";          I=0; J="Cambio"; K="Nodo Teco";    L="Pasante"; M=-58.455082; N=-34.558883; O="Saavedra";   P="Capital Norte" }
    @{ Row=5;  A="3839";      B="10/23/2024"; C="PICO 1511";                      D="13"; E="798390296"; F="NEW"; G="Pendiente";                      H="Poste inclinado"; I=1; J="Aplomo"; K="Sin equipos";  L="Poste";   M=-58.465596; N=-34.53627;  O="Saavedra";   P="Capital Norte" }
    @{ Row=6;  A="801645368"; B="12/13/2024"; C="San Blas 1809";                  D="11"; E="801645368"; F="NEW"; G="Pendiente";                      H="Picada";          I=0; J="Cambio"; K="Sin equipos";  L="Pasante"; M=-58.467767; N=-34.604588; O="Paternal";   P="Capital Norte" }
    @{ Row=7;  A="4595";      B="1/15/2025";  C="PAROISSIEN 1806";                D="13"; E="802747617"; F="NEW"; G="Pendiente";                      H="Aplomar";         I=1; J="Aplomo"; K="Sin equipos";  L="Terminal";M=-58.464172; N=-34.543845; O="Saavedra";   P="Capital Norte" }
    @{ Row=8;  A="4662";      B="1/21/2025";  C="ALTOLAGUIRRE 2397";              D="12"; E="802823938"; F="NEW"; G="Pendiente";                      H="Inclinada";       I=1; J="Aplomo"; K="Sin equipos";  L="Pasante"; M=-58.490766; N=-34.576987; O="Paternal";   P="Capital Norte" }
    @{ Row=69; A="7000";      B="8/20/2025";  C="SCALABRINI ORTIZ, RAUL AV. 2106";D="14"; E="809065867"; F="NEW"; G="Pendiente";                      H="Picada";          I=1; J="Cambio"; K="Sin equipos";  L="Pasante"; M=-58.420109; N=-34.58764;  O="Palermo";    P="Capital Sur" }
    @{ Row=70; A="7003";      B="8/20/2025";  C="MELIAN AV. 2576";                D="12"; E="809065874"; F="NEW"; G="Pendiente";                      H="Picada";          I=1; J="Cambio"; K="Sin equipos";  L="Pasante"; M=-58.471943; N=-34.564948; O="Colegiales"; P="Capital Norte" }
)

foreach ($r in $rowsData) {
    $row = $r.Row

    # Force the text-bearing columns to stay text even though several of
    # them (Caso/Comuna/OT numbers, dates) look numeric - otherwise Excel's
    # default "General" auto-detection would silently turn them into
    # numbers / dates on assignment.
    $ws.Range("A$row`:H$row").NumberFormat = "@"
    $ws.Range("J$row`:L$row").NumberFormat = "@"
    $ws.Range("O$row`:P$row").NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value  = $r.A
    $ws.Cells.Item($row, 2).Value  = $r.B
    $ws.Cells.Item($row, 3).Value  = $r.C
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = $r.E
    $ws.Cells.Item($row, 6).Value  = $r.F
    $ws.Cells.Item($row, 7).Value  = $r.G
    $ws.Cells.Item($row, 8).Value  = $r.H
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P

    # Clean up the temporary number-format back to the default style so we
    # don't leave a stray "@" text format applied to the cells.
    $ws.Range("A$row`:H$row").Style = "Normal"
    $ws.Range("J$row`:L$row").Style = "Normal"
    $ws.Range("O$row`:P$row").Style = "Normal"
}
